$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "36.317.40"
$ws.Range("E2").Value = "  +2.65%  "

$ws.Range("D3").Value = "2.007.92"
$ws.Range("E3").Value = "  +5.79%  "

$ws.Range("D4").Value = "'0.999"
$ws.Range("E4").Value = "  -0.13%  "

$ws.Range("D5").Value = "'245.23"
$ws.Range("E5").Value = "  -0.71%  "

$ws.Range("D6").Value = "'0.657"
$ws.Range("E6").Value = "  -4.81%  "

$ws.Range("D7").Value = "'0.998"
$ws.Range("E7").Value = "  -0.11%  "

$ws.Range("D8").Value = "'44.91"
$ws.Range("E8").Value = "  +4.54%  "

$ws.Range("D9").Value = "'0.361"
$ws.Range("E9").Value = "  +0.89%  "

$ws.Range("D10").Value = "'56.17"
$ws.Range("E10").Value = "  +0.77%  "

$ws.Range("D11").Value = "'0.0718"
$ws.Range("E11").Value = "  -3.79%  "

$ws.Range("E12").Value = "  +0.63%  "

$ws.Range("D13").Value = "'14.45"
$ws.Range("E13").Value = "  +3.37%  "

$ws.Range("D14").Value = "2.291.17"
$ws.Range("E14").Value = "  +5.67%  "

$ws.Range("D15").Value = "'0.795"
$ws.Range("E15").Value = "  +0.67%  "

$ws.Range("D16").Value = "1.999.06"
$ws.Range("E16").Value = "  +5.34%  "

$ws.Range("D17").Value = "'4.87"
$ws.Range("E17").Value = "  -2.43%  "

$ws.Range("D18").Value = "36.251.45"
$ws.Range("E18").Value = "  +2.53%  "

$ws.Range("D19").Value = "'70.72"
$ws.Range("E19").Value = "  -3.89%  "

$ws.Range("D20").Value = "0.0₃0811"
$ws.Range("E20").Value = "  -1.85%  "

$ws.Range("D21").Value = "'12.87"
$ws.Range("E21").Value = "  +0.03%  "

$ws.Range("D22").Value = "'233.31"
$ws.Range("E22").Value = "  -4.51%  "

$ws.Range("D23").Value = "'4.96"
$ws.Range("E23").Value = "  -5.59%  "

$ws.Range("D24").Value = "'1.00"
$ws.Range("E24").Value = "  +0.28%  "

$ws.Range("D25").Value = "'2.47"
$ws.Range("E25").Value = "  -7.09%  "

$ws.Range("D26").Value = "'161.42"
$ws.Range("E26").Value = "  -3.45%  "

$ws.Range("B27").Value = "EthereumClassic"
$ws.Range("C27").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D27").Value = "'19.73"
$ws.Range("E27").Value = "  +7.49%  "

$ws.Range("B28").Value = "PancakeSwap"
$ws.Range("C28").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D28").Value = "'1.97"
$ws.Range("E28").Value = "  -9.27%  "

$ws.Range("D29").Value = "'8.41"
$ws.Range("E29").Value = "  -1.38%  "

$ws.Range("D30").Value = "'0.122"
$ws.Range("E30").Value = "  -4.65%  "

$ws.Range("B31").Value = "Gas"
$ws.Range("C31").Value = "https://coinranking.com/coin/hfw0nnnLtSFc7+gas-gas"
$ws.Range("D31").Value = "'21.37"
$ws.Range("E31").Value = "  +58.98%  "

$ws.Range("B32").Value = "Filecoin"
$ws.Range("C32").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D32").Value = "'4.35"
$ws.Range("E32").Value = "  +0.10%  "

$ws.Range("D33").Value = "'0.0580"
$ws.Range("E33").Value = "  -3.28%  "

$ws.Range("D34").Value = "'0.998"
$ws.Range("E34").Value = "  -0.20%  "

$ws.Range("E35").Value = "  -0.13%  "

$ws.Range("B36").Value = "Kaspa"
$ws.Range("C36").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D36").Value = "'0.0841"
$ws.Range("E36").Value = "  +18.34%  "

$ws.Range("B37").Value = "InternetComputer(DFINITY)"
$ws.Range("C37").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D37").Value = "'4.00"
$ws.Range("E37").Value = "  -5.22%  "

$ws.Range("D38").Value = "'2.13"
$ws.Range("E38").Value = "  +8.66%  "

$ws.Range("D39").Value = "'0.835"
$ws.Range("E39").Value = "  -1.76%  "

$ws.Range("D40").Value = "'1.34"
$ws.Range("E40").Value = "  -8.60%  "

$ws.Range("D41").Value = "'96.77"
$ws.Range("E41").Value = "  -3.96%  "

$ws.Range("D42").Value = "'0.0214"
$ws.Range("E42").Value = "  -4.30%  "

$ws.Range("D43").Value = "'16.06"
$ws.Range("E43").Value = "  -6.27%  "

$ws.Range("D44").Value = "'1.07"
$ws.Range("E44").Value = "  -1.63%  "

$ws.Range("D45").Value = "'2.71"
$ws.Range("E45").Value = "  +12.04%  "

$ws.Range("D46").Value = "1.304.47"
$ws.Range("E46").Value = "  -2.36%  "

$ws.Range("D47").Value = "'0.0807"
$ws.Range("E47").Value = "  -0.41%  "

$ws.Range("E48").Value = "  +0.87%  "

$ws.Range("B49").Value = "RenderToken"
$ws.Range("C49").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D49").Value = "'2.20"
$ws.Range("E49").Value = "  -6.93%  "

$ws.Range("B50").Value = "RocketPoolETH"
$ws.Range("C50").Value = "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
$ws.Range("D50").Value = "2.180.05"
$ws.Range("E50").Value = "  +5.25%  "

$ws.Range("E51").Value = "  +10.47%  "
